$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.004727430580851
$ws.Cells.Item(2, 4).Value = 1.025507051647289
$ws.Cells.Item(2, 5).Value = 1.007762213487518
$ws.Cells.Item(2, 6).Value = 1.002867825820396
$ws.Cells.Item(2, 9).Value = 1.028047310659273
$ws.Cells.Item(2, 10).Value = 1.010014891386621
$ws.Cells.Item(2, 11).Value = 1.028332420894459
$ws.Cells.Item(2, 12).Value = 1.010640516912302
$ws.Cells.Item(2, 13).Value = 1.005761089072324
$ws.Cells.Item(2, 14).Value = 1.007465362557762

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.005822002272901
$ws.Cells.Item(3, 4).Value = 1.026019325364456
$ws.Cells.Item(3, 5).Value = 1.008694518224989
$ws.Cells.Item(3, 6).Value = 1.004617363073518
$ws.Cells.Item(3, 9).Value = 1.028114196249099
$ws.Cells.Item(3, 10).Value = 1.010738946473879
$ws.Cells.Item(3, 11).Value = 1.02865281339994
$ws.Cells.Item(3, 12).Value = 1.01137594559515
$ws.Cells.Item(3, 13).Value = 1.007310333110168
$ws.Cells.Item(3, 14).Value = 1.007712251451404

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.006530092566826
$ws.Cells.Item(4, 4).Value = 1.026344013162027
$ws.Cells.Item(4, 5).Value = 1.009297975763954
$ws.Cells.Item(4, 6).Value = 1.005748917304874
$ws.Cells.Item(4, 9).Value = 1.028152502439901
$ws.Cells.Item(4, 10).Value = 1.01120681070201
$ws.Cells.Item(4, 11).Value = 1.028852644648525
$ws.Cells.Item(4, 12).Value = 1.011851397160795
$ws.Cells.Item(4, 13).Value = 1.008311870664742
$ws.Cells.Item(4, 14).Value = 1.007871580600961

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.006827735276856
$ws.Cells.Item(5, 4).Value = 1.026478881440165
$ws.Cells.Item(5, 5).Value = 1.009551717083347
$ws.Cells.Item(5, 6).Value = 1.006224509550559
$ws.Cells.Item(5, 9).Value = 1.028167412619882
$ws.Cells.Item(5, 10).Value = 1.011403347421962
$ws.Cells.Item(5, 11).Value = 1.028934856412712
$ws.Cells.Item(5, 12).Value = 1.012051177842443
$ws.Cells.Item(5, 13).Value = 1.008732703681971
$ws.Cells.Item(5, 14).Value = 1.007938461346076

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.006877708613519
$ws.Cells.Item(6, 4).Value = 1.026501430639359
$ws.Cells.Item(6, 5).Value = 1.00959432420678
$ws.Cells.Item(6, 6).Value = 1.006304357228488
$ws.Cells.Item(6, 9).Value = 1.0281698460087
$ws.Cells.Item(6, 10).Value = 1.011436337855459
$ws.Cells.Item(6, 11).Value = 1.028948554567337
$ws.Cells.Item(6, 12).Value = 1.012084716123469
$ws.Cells.Item(6, 13).Value = 1.008803351132689
$ws.Cells.Item(6, 14).Value = 1.007949684992606

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.006534069831341
$ws.Cells.Item(7, 4).Value = 1.026345821692922
$ws.Cells.Item(7, 5).Value = 1.009301366078558
$ws.Cells.Item(7, 6).Value = 1.005755272624287
$ws.Cells.Item(7, 9).Value = 1.028152706365405
$ws.Cells.Item(7, 10).Value = 1.01120943743596
$ws.Cells.Item(7, 11).Value = 1.028853750236068
$ws.Cells.Item(7, 12).Value = 1.011854067027955
$ws.Cells.Item(7, 13).Value = 1.008317494683953
$ws.Cells.Item(7, 14).Value = 1.007872474662037

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.005097383039582
$ws.Cells.Item(8, 4).Value = 1.025681579129091
$ws.Cells.Item(8, 5).Value = 1.008077250859506
$ws.Cells.Item(8, 6).Value = 1.003459202297312
$ws.Cells.Item(8, 9).Value = 1.02807094286734
$ws.Cells.Item(8, 10).Value = 1.010259724350782
$ws.Cells.Item(8, 11).Value = 1.028442245416833
$ws.Cells.Item(8, 12).Value = 1.010889146272544
$ws.Cells.Item(8, 13).Value = 1.006284860424486
$ws.Cells.Item(8, 14).Value = 1.00754888794167

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.002564339795591
$ws.Cells.Item(9, 4).Value = 1.024459342914767
$ws.Cells.Item(9, 5).Value = 1.005921629741186
$ws.Cells.Item(9, 6).Value = 0.9994088644037524
$ws.Cells.Item(9, 9).Value = 1.027888906352644
$ws.Cells.Item(9, 10).Value = 1.008581170752679
$ws.Cells.Item(9, 11).Value = 1.027660033183786
$ws.Cells.Item(9, 12).Value = 1.009185553564632
$ws.Cells.Item(9, 13).Value = 1.002695617146474
$ws.Cells.Item(9, 14).Value = 1.00697541805312

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.000874535888566
$ws.Cells.Item(10, 4).Value = 1.023610008354521
$ws.Cells.Item(10, 5).Value = 1.004485421829613
$ws.Cells.Item(10, 6).Value = 0.9967050907846262
$ws.Cells.Item(10, 9).Value = 1.027742192411911
$ws.Cells.Item(10, 10).Value = 1.007458639509106
$ws.Cells.Item(10, 11).Value = 1.027100466945094
$ws.Cells.Item(10, 12).Value = 1.008047530887223
$ws.Cells.Item(10, 13).Value = 1.000297222101511
$ws.Cells.Item(10, 14).Value = 1.006590881253872

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.000142534010965
$ws.Cells.Item(11, 4).Value = 1.02323410637176
$ws.Cells.Item(11, 5).Value = 1.003863715011103
$ws.Cells.Item(11, 6).Value = 0.995533349279506
$ws.Cells.Item(11, 9).Value = 1.027672681232815
$ws.Cells.Item(11, 10).Value = 1.006971717307854
$ws.Cells.Item(11, 11).Value = 1.026849189472056
$ws.Cells.Item(11, 12).Value = 1.007554189785336
$ws.Cells.Item(11, 13).Value = 0.999257251006217
$ws.Cells.Item(11, 14).Value = 1.006423839134035

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 0.99987058688519
$ws.Cells.Item(12, 4).Value = 1.023093261465912
$ws.Cells.Item(12, 5).Value = 1.003632810757945
$ws.Cells.Item(12, 6).Value = 0.9950979518761309
$ws.Cells.Item(12, 9).Value = 1.027645965106361
$ws.Cells.Item(12, 10).Value = 1.006790721502815
$ws.Cells.Item(12, 11).Value = 1.026754508480939
$ws.Cells.Item(12, 12).Value = 1.007370853485298
$ws.Cells.Item(12, 13).Value = 0.9988707312407346
$ws.Cells.Item(12, 14).Value = 1.006361711295728

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 0.9999289227360615
$ws.Cells.Item(13, 4).Value = 1.02312352824074
$ws.Cells.Item(13, 5).Value = 1.00368233938652
$ws.Cells.Item(13, 6).Value = 0.9951913535466643
$ws.Cells.Item(13, 9).Value = 1.027651736342748
$ws.Cells.Item(13, 10).Value = 1.006829551690842
$ws.Cells.Item(13, 11).Value = 1.02677487869296
$ws.Cells.Item(13, 12).Value = 1.007410183731069
$ws.Cells.Item(13, 13).Value = 0.9989536515352102
$ws.Cells.Item(13, 14).Value = 1.006375041600458

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.000120055795136
$ws.Cells.Item(14, 4).Value = 1.023222488916466
$ws.Cells.Item(14, 5).Value = 1.003844627887068
$ws.Cells.Item(14, 6).Value = 0.9954973625354055
$ws.Cells.Item(14, 9).Value = 1.027670491146611
$ws.Cells.Item(14, 10).Value = 1.006956758827764
$ws.Cells.Item(14, 11).Value = 1.026841390534216
$ws.Cells.Item(14, 12).Value = 1.007539036944252
$ws.Cells.Item(14, 13).Value = 0.9992253058768764
$ws.Cells.Item(14, 14).Value = 1.006418705282296

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.000237812589679
$ws.Cells.Item(15, 4).Value = 1.023283300548361
$ws.Cells.Item(15, 5).Value = 1.003944622492043
$ws.Cells.Item(15, 6).Value = 0.9956858831607941
$ws.Cells.Item(15, 9).Value = 1.027681927861248
$ws.Cells.Item(15, 10).Value = 1.007035117889274
$ws.Cells.Item(15, 11).Value = 1.026882192518695
$ws.Cells.Item(15, 12).Value = 1.007618416018253
$ws.Cells.Item(15, 13).Value = 0.9993926505345242
$ws.Cells.Item(15, 14).Value = 1.006445597173323

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.000923109737314
$ws.Cells.Item(16, 4).Value = 1.023634784689578
$ws.Cells.Item(16, 5).Value = 1.004526686116763
$ws.Cells.Item(16, 6).Value = 0.9967828335567679
$ws.Cells.Item(16, 9).Value = 1.027746679784026
$ws.Cells.Item(16, 10).Value = 1.007490936651406
$ws.Cells.Item(16, 11).Value = 1.02711695452225
$ws.Cells.Item(16, 12).Value = 1.008080260149639
$ws.Cells.Item(16, 13).Value = 1.000366210124534
$ws.Cells.Item(16, 14).Value = 1.006601955966833

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.001352894509748
$ws.Cells.Item(17, 4).Value = 1.023853086440131
$ws.Cells.Item(17, 5).Value = 1.00489184635577
$ws.Cells.Item(17, 6).Value = 0.997470648322167
$ws.Cells.Item(17, 9).Value = 1.027785696793537
$ws.Cells.Item(17, 10).Value = 1.007776628254542
$ws.Cells.Item(17, 11).Value = 1.027261812832631
$ws.Cells.Item(17, 12).Value = 1.00836980892789
$ws.Cells.Item(17, 13).Value = 1.000976502807035
$ws.Cells.Item(17, 14).Value = 1.006699892041016

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.00160355141928
$ws.Cells.Item(18, 4).Value = 1.02397963304905
$ws.Cells.Item(18, 5).Value = 1.005104855702431
$ws.Cells.Item(18, 6).Value = 0.9978717443451044
$ws.Cells.Item(18, 9).Value = 1.027807877329898
$ws.Cells.Item(18, 10).Value = 1.007943184624868
$ws.Cells.Item(18, 11).Value = 1.027345439387844
$ws.Cells.Item(18, 12).Value = 1.008538642982987
$ws.Cells.Item(18, 13).Value = 1.001332336979295
$ws.Cells.Item(18, 14).Value = 1.006756964900417

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.001689014075675
$ws.Cells.Item(19, 4).Value = 1.024022648879114
$ws.Cells.Item(19, 5).Value = 1.005177489416966
$ws.Cells.Item(19, 6).Value = 0.9980084920162302
$ws.Cells.Item(19, 9).Value = 1.027815342305294
$ws.Cells.Item(19, 10).Value = 1.007999962080195
$ws.Cells.Item(19, 11).Value = 1.027373806722037
$ws.Cells.Item(19, 12).Value = 1.008596201784786
$ws.Cells.Item(19, 13).Value = 1.001453643932613
$ws.Cells.Item(19, 14).Value = 1.006776416535688

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.001306785723031
$ws.Cells.Item(20, 4).Value = 1.023829745907683
$ws.Cells.Item(20, 5).Value = 1.004852666315862
$ws.Cells.Item(20, 6).Value = 0.9973968621890754
$ws.Cells.Item(20, 9).Value = 1.027781570346936
$ws.Cells.Item(20, 10).Value = 1.007745984820907
$ws.Cells.Item(20, 11).Value = 1.027246360540018
$ws.Cells.Item(20, 12).Value = 1.008338748771612
$ws.Cells.Item(20, 13).Value = 1.000911038655298
$ws.Cells.Item(20, 14).Value = 1.006689389768212

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.000063773251871
$ws.Cells.Item(21, 4).Value = 1.023193381061208
$ws.Cells.Item(21, 5).Value = 1.003796837283646
$ws.Cells.Item(21, 6).Value = 0.9954072550388525
$ws.Cells.Item(21, 9).Value = 1.027664993056495
$ws.Cells.Item(21, 10).Value = 1.006919303136536
$ws.Cells.Item(21, 11).Value = 1.026821841561943
$ws.Cells.Item(21, 12).Value = 1.007501095313476
$ws.Cells.Item(21, 13).Value = 0.9991453168265775
$ws.Cells.Item(21, 14).Value = 1.006405849655614

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 0.9992819567113377
$ws.Cells.Item(22, 4).Value = 1.022786227006253
$ws.Cells.Item(22, 5).Value = 1.003133141372291
$ws.Cells.Item(22, 6).Value = 0.9941553759467042
$ws.Cells.Item(22, 9).Value = 1.027586509990163
$ws.Cells.Item(22, 10).Value = 1.006398775255053
$ws.Cells.Item(22, 11).Value = 1.026547147649343
$ws.Cells.Item(22, 12).Value = 1.006973922460185
$ws.Cells.Item(22, 13).Value = 0.9980338120025594
$ws.Cells.Item(22, 14).Value = 1.006227108057501

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 0.999696440536912
$ws.Cells.Item(23, 4).Value = 1.023002733729557
$ws.Cells.Item(23, 5).Value = 1.003484965927849
$ws.Cells.Item(23, 6).Value = 0.994819113223859
$ws.Cells.Item(23, 9).Value = 1.027628606217047
$ws.Cells.Item(23, 10).Value = 1.006674789739233
$ws.Cells.Item(23, 11).Value = 1.026693504457306
$ws.Cells.Item(23, 12).Value = 1.00725343543061
$ws.Cells.Item(23, 13).Value = 0.9986231708368029
$ws.Cells.Item(23, 14).Value = 1.006321906973259

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.001327620389373
$ws.Cells.Item(24, 4).Value = 1.023840294915894
$ws.Cells.Item(24, 5).Value = 1.004870370035164
$ws.Cells.Item(24, 6).Value = 0.9974302032580507
$ws.Cells.Item(24, 9).Value = 1.027783436694908
$ws.Cells.Item(24, 10).Value = 1.007759831525004
$ws.Cells.Item(24, 11).Value = 1.027253345444725
$ws.Cells.Item(24, 12).Value = 1.008352783688695
$ws.Cells.Item(24, 13).Value = 1.000940619516847
$ws.Cells.Item(24, 14).Value = 1.006694135452648

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.003219376118636
$ws.Cells.Item(25, 4).Value = 1.024781422877932
$ws.Cells.Item(25, 5).Value = 1.006478749260444
$ws.Cells.Item(25, 6).Value = 1.00045655425395
$ws.Cells.Item(25, 9).Value = 1.027940447854658
$ws.Cells.Item(25, 10).Value = 1.009015725207216
$ws.Cells.Item(25, 11).Value = 1.027868987877461
$ws.Cells.Item(25, 12).Value = 1.009626370872034
$ws.Cells.Item(25, 13).Value = 1.003624464200099
$ws.Cells.Item(25, 14).Value = 1.007124063766357

Write-Host "Applied vm_pu.xlsx updates for 380 kV case"
